# Automatische test-sync: 2025-07-22 18:01:50
# Appends a new "Openingstijden / Locatie" test-mail row to the Logs sheet
# and bumps the matching count on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 22

$logs.Cells.Item($row, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Cells.Item($row, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($row, 5).Value = "Beste klant,`r`n`r`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`r`n`r`nMet vriendelijke groet,`r`n`r`n[Bedrijfsnaam]"
$logs.Cells.Item($row, 6).Value = "2025-07-22 18:01:47"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"
$logs.Cells.Item($row, 10).Value = "Ja"

# The new cell's text contains embedded line breaks, which makes the engine
# auto-grow the row height (like Excel does on entry). Re-run AutoFit so the
# row keeps its normal/default height instead of persisting a custom one.
$logs.Rows.Item($row).AutoFit()

# Extend the existing conditional formatting ranges (D/G/H/I/J) so row 22 is
# covered as well, same as the rest of the table.
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))
$logs.Range("H2:H21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H22"))
$logs.Range("I2:I21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I22"))
$logs.Range("J2:J21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J22"))

# Bump the "Openingstijden / Locatie" count on the dashboard summary sheet.
$dashboard.Cells.Item(4, 2).Value = 4
